# Update kanji flashcard deck:
#  - Slides 1-9: swap in new kanji/definition text, update page reference 67-68 -> 65-66
#  - Add a new Slide 10 (duplicated from Slide 9's layout/formatting) for 害虫
#  - Give the new slide a notes page (slide-number field auto-populates to 10)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the nine existing slides (title / definition / page-range shapes)
# ---------------------------------------------------------------------------

$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = '幅広い'
$s1.Shapes.Item(3).TextFrame.TextRange.Text = 'extensive, wide, broad...'
$s1.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = '増幅'
$s2.Shapes.Item(3).TextFrame.TextRange.Text = 'amplification (elec.) | magnification, amplification, making larger...'
$s2.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = '支える'
$s3.Shapes.Item(3).TextFrame.TextRange.Text = 'to support, to prop, to sustain, to underlay, to hold up, to defend | to hold at bay, to stem, to check...'
$s3.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = '支持'
$s4.Shapes.Item(3).TextFrame.TextRange.Text = 'support, backing, endorsement, approval | propping up, holding up, support...'
$s4.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = '反対'
$s5.Shapes.Item(3).TextFrame.TextRange.Text = 'opposition, resistance, antagonism, hostility, objection, dissent | reverse, opposite, inverse, contrary...'
$s5.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = '論文'
$s6.Shapes.Item(3).TextFrame.TextRange.Text = 'thesis, essay, treatise, paper, article...'
$s6.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = '否定'
$s7.Shapes.Item(3).TextFrame.TextRange.Text = 'denial, negation, repudiation, disavowal | negation | NOT operation...'
$s7.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = '野菜'
$s8.Shapes.Item(3).TextFrame.TextRange.Text = 'vegetable...'
$s8.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = '穀物'
$s9.Shapes.Item(3).TextFrame.TextRange.Text = 'grain, cereal, corn...'
$s9.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

# ---------------------------------------------------------------------------
# 2. Add Slide 10 - duplicate Slide 9 so the new slide inherits the exact
#    textbox layout/formatting used throughout the deck, then retarget text.
# ---------------------------------------------------------------------------

$dup = $p.Slides.Item(9).Duplicate()
$s10 = $p.Slides.Item(10)
$s10.Name = 'Slide 10'

$s10.Shapes.Item(1).TextFrame.TextRange.Text = '害虫'
$s10.Shapes.Item(3).TextFrame.TextRange.Text = 'harmful insect, noxious insect, vermin, pest...'
$s10.Shapes.Item(4).TextFrame.TextRange.Text = '65-66'

# ---------------------------------------------------------------------------
# 3. Materialize a notes page for the new slide (slide-number field shows 10)
# ---------------------------------------------------------------------------

$notes10 = $s10.NotesPage
$notesBody = $notes10.Shapes.Placeholders.Item(2)
$notesBody.TextFrame.TextRange.Text = ''
